# Update column F (dSF) values for rows 2-14, 16, 17, 19 per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 1
    6  = -2
    7  = -2
    8  = -1
    9  = 2
    10 = 2
    11 = -2
    12 = 2
    13 = 1
    14 = 2
    16 = 1
    17 = 6
    19 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
